# Auto-generated PowerShell COM-interop script
# Applies numeric cell updates to match the target market-price refresh diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 100
$ws.Range("K5").Value = 100
$ws.Range("M5").Value = 15
$ws.Range("H63").Value = 32325
$ws.Range("J63").Value = 32325
$ws.Range("L63").Value = 32325
$ws.Range("N63").Value = -33573
$ws.Range("H66").Value = 32325
$ws.Range("J66").Value = 32325
$ws.Range("L66").Value = 96975
$ws.Range("N66").Value = -103215
$ws.Range("H81").Value = 19877.3
$ws.Range("I81").Value = 10249
$ws.Range("J81").Value = 22284.375
$ws.Range("K81").Value = 10249
$ws.Range("L81").Value = 22284.375
$ws.Range("M81").Value = -9251
$ws.Range("N81").Value = -24280.375
$ws.Range("H84").Value = 19877.3
$ws.Range("I84").Value = 10249
$ws.Range("J84").Value = 22284.375
$ws.Range("K84").Value = 30747
$ws.Range("L84").Value = 66853.125
$ws.Range("M84").Value = -25755
$ws.Range("N84").Value = -76837.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2376.2307
$ws.Range("I122").Value = 1494.6364
$ws.Range("J122").Value = 7225
$ws.Range("K122").Value = 4483.9092
$ws.Range("L122").Value = 21675
$ws.Range("M122").Value = -2033.9092
$ws.Range("N122").Value = -26575

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 3720.8333
$ws.Range("J62").Value = 3880.8823
$ws.Range("L62").Value = 11642.6469
$ws.Range("N62").Value = -13014.6469
$ws.Range("H65").Value = 3720.8333
$ws.Range("J65").Value = 3880.8823
$ws.Range("L65").Value = 34927.9407
$ws.Range("N65").Value = -41791.9407
$ws.Range("H68").Value = 2263.3242
$ws.Range("I68").Value = 811.25
$ws.Range("J68").Value = 2663.8965
$ws.Range("K68").Value = 2433.75
$ws.Range("L68").Value = 7991.689499999999
$ws.Range("M68").Value = -1622.75
$ws.Range("N68").Value = -9613.6895
$ws.Range("H69").Value = 14212.25
$ws.Range("J69").Value = 18666.5
$ws.Range("L69").Value = 55999.5
$ws.Range("N69").Value = -57621.5
$ws.Range("H70").Value = 4697
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 4697
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 14091
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -14721
$ws.Range("H71").Value = 2263.3242
$ws.Range("I71").Value = 811.25
$ws.Range("J71").Value = 2663.8965
$ws.Range("K71").Value = 7301.25
$ws.Range("L71").Value = 23975.0685
$ws.Range("M71").Value = -3245.25
$ws.Range("N71").Value = -32087.0685
$ws.Range("H72").Value = 14212.25
$ws.Range("J72").Value = 18666.5
$ws.Range("L72").Value = 167998.5
$ws.Range("N72").Value = -176110.5
$ws.Range("H73").Value = 4697
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 4697
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 14091
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -16275
$ws.Range("H74").Value = 23000
$ws.Range("J74").Value = 28000
$ws.Range("L74").Value = 84000
$ws.Range("N74").Value = -86122
$ws.Range("H75").Value = 2974.375
$ws.Range("I75").Value = 940
$ws.Range("J75").Value = 3652.5
$ws.Range("K75").Value = 2820
$ws.Range("L75").Value = 10957.5
$ws.Range("M75").Value = -1822
$ws.Range("N75").Value = -12953.5
$ws.Range("H76").Value = 3875
$ws.Range("I76").Value = 3500
$ws.Range("K76").Value = 10500
$ws.Range("M76").Value = -10117
$ws.Range("H77").Value = 23000
$ws.Range("J77").Value = 28000
$ws.Range("L77").Value = 252000
$ws.Range("N77").Value = -262608
$ws.Range("H78").Value = 2974.375
$ws.Range("I78").Value = 940
$ws.Range("J78").Value = 3652.5
$ws.Range("K78").Value = 8460
$ws.Range("L78").Value = 32872.5
$ws.Range("M78").Value = -3468
$ws.Range("N78").Value = -42856.5
$ws.Range("H79").Value = 3875
$ws.Range("I79").Value = 3500
$ws.Range("K79").Value = 10500
$ws.Range("M79").Value = -9174
$ws.Range("H80").Value = 16903
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H81").Value = 36000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 36000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 108000
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -110246
$ws.Range("H82").Value = 4000
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H83").Value = 16903
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H84").Value = 36000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 36000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 324000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -335232
$ws.Range("H85").Value = 4000
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H86").Value = 1800.5
$ws.Range("J86").Value = 1800.5
$ws.Range("L86").Value = 5401.5
$ws.Range("N86").Value = -7773.5
$ws.Range("H87").Value = 15000
$ws.Range("J87").Value = 15000
$ws.Range("L87").Value = 45000
$ws.Range("N87").Value = -47496
$ws.Range("H88").Value = 6002
$ws.Range("J88").Value = 6002
$ws.Range("L88").Value = 18006
$ws.Range("N88").Value = -18862
$ws.Range("H89").Value = 1800.5
$ws.Range("J89").Value = 1800.5
$ws.Range("L89").Value = 16204.5
$ws.Range("N89").Value = -28060.5
$ws.Range("H90").Value = 15000
$ws.Range("J90").Value = 15000
$ws.Range("L90").Value = 135000
$ws.Range("N90").Value = -147480
$ws.Range("H91").Value = 6002
$ws.Range("J91").Value = 6002
$ws.Range("L91").Value = 18006
$ws.Range("N91").Value = -20970

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4886.2144
$ws.Range("I122").Value = 3214.1428
$ws.Range("J122").Value = 6558.2856
$ws.Range("K122").Value = 9642.428400000001
$ws.Range("L122").Value = 19674.8568
$ws.Range("M122").Value = -7192.428400000001
$ws.Range("N122").Value = -24574.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2715
$ws.Range("I40").Value = 1520
$ws.Range("J40").Value = 2885.7144
$ws.Range("K40").Value = 1520
$ws.Range("L40").Value = 2885.7144
$ws.Range("M40").Value = -1384
$ws.Range("N40").Value = -3157.7144

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 36397.832
$ws.Range("J46").Value = 36397.832
$ws.Range("L46").Value = 36397.832
$ws.Range("N46").Value = -36859.832
$ws.Range("H80").Value = 30000
$ws.Range("J80").Value = 30000
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -31996
$ws.Range("H83").Value = 30000
$ws.Range("J83").Value = 30000
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -99984
$ws.Range("H132").Value = 6097.3613
$ws.Range("I132").Value = 2692.577
$ws.Range("J132").Value = 14949.8
$ws.Range("K132").Value = 8077.731000000001
$ws.Range("L132").Value = 44849.39999999999
$ws.Range("M132").Value = -5547.731000000001
$ws.Range("N132").Value = -49909.39999999999
$ws.Range("H134").Value = 36397.832
$ws.Range("J134").Value = 36397.832
$ws.Range("L134").Value = 109193.496
$ws.Range("N134").Value = -114263.496
